$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Datos" to "diagnosticos_eventos"
$ws.Name = "diagnosticos_eventos"

# Remove the custom header styling (bold white font on blue fill, centered)
# from the header row A1:G1 so the cells fall back to the default style.
$ws.Range("A1:G1").ClearFormats()

# Drop the hard-coded 20-character column widths on A:G and restore the
# sheet default width. Deleting the column range clears its stored
# dimension/width metadata; re-typing the header values below puts the
# (now unstyled, default-width) columns back exactly as they were.
$ws.Columns("A:G").Delete()

$ws.Range("A1").Value = "animal_codigo"
$ws.Range("B1").Value = "fecha"
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "diagnostico_detalle"
$ws.Range("E1").Value = "severidad"
$ws.Range("F1").Value = "estado"
$ws.Range("G1").Value = "observaciones"
